$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the format/style of the existing
# header cells (e.g. G1 "sum") by copying formats over.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add data values for the new Save column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
